# Generate Report for Handoff
#
# The localization status report is regenerated for handoff: the
# "In Translation" status becomes "Ready for handoff" on every sheet that
# carries it (Overview!E2/F2, zh-cn!C2, de-de!C2), and the associated
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps are
# bumped to the new handoff-generation time. Because "Ready for handoff" is
# wider than "In Translation", the Status columns are re-fitted to the new
# text width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Latest HO Xliff Generate Date / Latest Handoff Datetime ---
$wsOverview.Range("G2").Value = "2016-08-27 00:39:57"
$wsDeDe.Range("H2").Value = "2016-08-27 00:39:57"
$wsZhCn.Range("H2").Value = "2016-08-27 00:39:52"

# --- Re-fit the Status columns now that the text is wider ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
